# UCLA / Irvine / Modified Colleges
# Adds "Class" (E), "Full Time" (F), "Certificate Program" (G) and
# "Region" (H) data to every data row (2-50) of the UCIrvine sheet, and
# updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 50

# --- Column E ("Class"): same shared string + style (s="1") on every
#     data row, copied from the style already used by column A/B/C. ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "Class"
}
$ws.Cells.Item(2, 1).Copy()
$ws.Range("E2:E$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Columns F ("Full Time") and G ("Certificate Program"): plain
#     numbers, no special style. Rows 2-8 are Full Time=1 Certificate=1;
#     rows 9-50 are Full Time=1 Certificate=0. ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
    if ($r -le 8) {
        $ws.Cells.Item($r, 7).Value = 1
    } else {
        $ws.Cells.Item($r, 7).Value = 0
    }
}

# --- Column H ("Region"): "Southern" for every data row (row 2 already
#     had it before this edit). Reuses the existing shared string, no
#     style change needed. ---
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "Southern"
}

# --- Update the sheet's active selection to match the latest edit. ---
$ws.Range("H53").Select()
